$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AD1").Font.Name = "Comic Sans MS"
